# ADD: add language info of column when term_expansion and revise
#
# Target sheet: "amazon" (second tab in the workbook).
# 1. Add a new "lang" column (M) with "English" on the rows that already
#    carry a term_expansion value (term_expansion style J11/J4/J5/J15/J16).
# 2. Rewrite the "term_expansion" placeholder text in column J for the
#    brand/compatible/scenario rows (19-21) from the old
#    "#VIRTUAL_IN(...)" form to the new "VIRTUAL_IN(....%{value}%)" form.
# 3. Update the current selection to E24 (matches the saved sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("amazon")

# --- New "lang" column header + values -------------------------------
# Write the English values before the header so the shared-strings table
# grows in the same order as the authored workbook (English, lang, ...).
$ws.Range("M4").Value = "English"
$ws.Range("M1").Value = "lang"

$ws.Range("M5").Value = "English"
$ws.Range("M11").Value = "English"
$ws.Range("M15").Value = "English"
$ws.Range("M16").Value = "English"

# --- Revised term_expansion placeholders ------------------------------
$ws.Range("J21").Value = "VIRTUAL_IN(category.%{value}%)"
$ws.Range("J19").Value = "VIRTUAL_IN(product_name.%{value}%)"
$ws.Range("J20").Value = "VIRTUAL_IN(product_name.%{value}%)"

# --- Update the saved selection ---------------------------------------
$ws.Range("E24").Select()
